# Golden R01 Behavior Master Key - Run 5/6 update
# Adds 23 new rows (Run 6 behavior) to the master table, extending it
# from A1:X105 to A1:X128, and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ A=73; B=105; C=6; E=45660; F="Male"; G="c57"; I=9; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=74; B=106; C=6; E=45660; F="Male"; G="c57"; I=10; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=75; B=107; C=6; E=45660; F="Male"; G="c57"; I=11; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=76; B=108; C=6; E=45660; F="Male"; G="c57"; I=12; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=77; B=109; C=6; E=45660; F="Male"; G="c57"; I=13; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=78; B=110; C=6; E=45660; F="Male"; G="c57"; I=14; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=79; B=111; C=6; E=45660; F="Female"; G="c57"; I=15; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=80; B=112; C=6; E=45660; F="Female"; G="c57"; I=16; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=81; B=113; C=6; E=45660; F="Female"; G="c57"; I=17; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=82; B=114; C=6; E=45660; F="Female"; G="c57"; I=18; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=83; B=115; C=6; E=45660; F="Female"; G="c57"; I=19; J="Morning"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=84; B=116; C=6; E=45660; F="Female"; G="c57"; I=20; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=85; B=117; C=6; E=45660; F="Male"; G="c57"; I=9; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=86; B=118; C=6; E=45660; F="Male"; G="c57"; I=10; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=87; B=119; C=6; E=45660; F="Male"; G="c57"; I=11; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=88; B=120; C=6; E=45660; F="Male"; G="c57"; I=12; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=89; B=121; C=6; E=45660; F="Male"; G="c57"; I=13; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=90; B=122; C=6; E=45660; F="Female"; G="c57"; I=15; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=91; B=123; C=6; E=45660; F="Female"; G="c57"; I=16; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=92; B=124; C=6; E=45660; F="Female"; G="c57"; I=17; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=93; B=125; C=6; E=45660; F="Female"; G="c57"; I=18; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=94; B=126; C=6; E=45660; F="Female"; G="c57"; I=19; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false },
  @{ A=95; B=127; C=6; E=45660; F="Female"; G="c57"; I=20; J="Afternoon"; K=$true; L=$true; M=$true; N=$false; O=$false }
)

# Expand the table ("Table1") by the number of new rows first, so the
# table definition / autofilter range grows along with the data.
$lo = $ws.ListObjects.Item("Table1")
for ($i = 0; $i -lt $rows.Count; $i++) {
    $lo.ListRows.Add() | Out-Null
}

$startRow = 106
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A   # TagNumber
    $ws.Cells.Item($r, 2).Value = $row.B   # ID
    $ws.Cells.Item($r, 3).Value = $row.C   # Run
    $eCell = $ws.Cells.Item($r, 5)         # DOB
    $eCell.Value = $row.E
    $eCell.NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 6).Value = $row.F   # Sex
    $ws.Cells.Item($r, 7).Value = $row.G   # Strain
    $ws.Cells.Item($r, 9).Value = $row.I   # Chamber
    $ws.Cells.Item($r, 10).Value = $row.J  # TimeOfBehavior
    $ws.Cells.Item($r, 11).Value = $row.K  # SelfAdministration
    $ws.Cells.Item($r, 12).Value = $row.L  # Extinction
    $ws.Cells.Item($r, 13).Value = $row.M  # Reinstatement
    $ws.Cells.Item($r, 14).Value = $row.N  # BehavioralEconomics
    $ws.Cells.Item($r, 15).Value = $row.O  # ProgressiveRatio
    $r++
}

# Restore the view: scroll/select as it was left after the edit.
$ws.Range("D132").Select()
